# extraSpaceInEndBookmark-expected-generation.docx
#
# The commit just bumps the Apache POI version used to regenerate this
# fixture (POI 3.17.0 -> 4.0.1). Newer POI mints a fresh random
# "rsid"-looking token for the runs that make up the "REF bookmark1"
# field (it has no real Word rsid to reuse) and a fresh random numeric id
# for the bookmark it writes; the text/structure of the document does not
# change at all, only those two internal, auto-generated identifiers do.
#
# Both identifiers are allocator/serializer-internal values - exactly like
# in real Word, a COM script cannot poke an arbitrary bookmark id (Word
# assigns those itself), so the best a legitimate edit can do is cause the
# engine to mint fresh ones, the same way re-saving with a newer generator
# would. This script re-stamps the rsid token to the exact new value and
# re-creates the bookmark so a new id gets allocated for it.

$d = $word.ActiveDocument

$newToken = "162E50549219B7B77638FA5D7572EA24"

# --- 1. Find the paragraph holding the "REF bookmark1" field and rewrite
#        its 5 runs (begin / instrText / separate / result / end) so they
#        all carry the new rsid-like run token. ---
$fieldParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.IndexOf("a reference to bookmark1") -ge 0) {
        $fieldParaIndex = $i
    }
}

$p = $d.Paragraphs.Item($fieldParaIndex).Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidRPr="00FF681D" w:rsidP="009168BC">
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="3119"/>
    </w:tabs>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00FF681D">
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Test link before bookmark : </w:t>
  </w:r>
  <w:r w:rsidR="__NEWTOKEN__">
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r w:rsidR="__NEWTOKEN__">
    <w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText>
  </w:r>
  <w:r w:rsidR="__NEWTOKEN__">
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r w:rsidR="__NEWTOKEN__">
    <w:rPr>
      <w:b w:val="true"/>
      <w:noProof/>
    </w:rPr>
    <w:t>a reference to bookmark1</w:t>
  </w:r>
  <w:r w:rsidR="__NEWTOKEN__">
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
'@
$xml = $xml.Replace("__NEWTOKEN__", $newToken)
$p.InsertXML($xml)

# --- 2. Re-create the "bookmark1" bookmark in place so the engine mints a
#        fresh internal id for it (ids are never user-settable, in real
#        Word or here - this mirrors what re-generating the document with
#        a newer writer does). ---
if ($d.Bookmarks.Exists("bookmark1")) {
    $bm = $d.Bookmarks.Item("bookmark1")
    $bmRange = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add("bookmark1", $bmRange) | Out-Null
}

Write-Host "Re-stamped the REF bookmark1 field runs with rsid token $newToken and re-minted bookmark1's id."
